$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")

$ws0.Range("B2").Value = -0.3083888998137556
$ws0.Range("C2").Value = -0.3266766542974147

$ws0.Range("B3").Value = -0.6934534984498852
$ws0.Range("C3").Value = 0.3008124411097155

$ws0.Range("B4").Value = -1.298486295815009
$ws0.Range("C4").Value = -1.211234106537655

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")

$ws1.Range("B2").Value = -0.665847225197116
$ws1.Range("C2").Value = 0.07973425602387073

$ws1.Range("B3").Value = -1.159976077872662
$ws1.Range("C3").Value = -0.7025142724985503

$ws1.Range("B4").Value = -0.5629205337645407
$ws1.Range("C4").Value = 0.253028877470333
